$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "68.192.79"
$ws.Range("E2").Value = "  -0.91%  "
Set-TextValue $ws.Range("D3") "3.888.77"
$ws.Range("E3").Value = "  -0.96%  "
$ws.Range("E4").Value = "  -0.05%  "
Set-TextValue $ws.Range("D5") "600.25"
$ws.Range("E5").Value = "  -0.46%  "
Set-TextValue $ws.Range("D6") "171.07"
$ws.Range("E6").Value = "  +1.95%  "
Set-TextValue $ws.Range("D7") "3.886.75"
$ws.Range("E7").Value = "  -1.04%  "
$ws.Range("E8").Value = "  +0.14%  "
Set-TextValue $ws.Range("D9") "0.531"
$ws.Range("E9").Value = "  -0.40%  "
$ws.Range("E10").Value = "  -3.02%  "
Set-TextValue $ws.Range("D11") "6.41"
$ws.Range("E11").Value = "  -0.88%  "
Set-TextValue $ws.Range("D12") "0.458"
$ws.Range("E12").Value = "  -1.27%  "
Set-TextValue $ws.Range("D13") "0.0000259"
$ws.Range("E13").Value = "  +0.90%  "
Set-TextValue $ws.Range("D14") "37.12"
$ws.Range("E14").Value = "  -1.36%  "
Set-TextValue $ws.Range("D15") "4.542.71"
$ws.Range("E15").Value = "  -0.84%  "
Set-TextValue $ws.Range("D16") "3.900.58"
$ws.Range("E16").Value = "  -1.31%  "
Set-TextValue $ws.Range("D17") "68.338.76"
$ws.Range("E17").Value = "  -0.80%  "
Set-TextValue $ws.Range("D18") "18.12"
$ws.Range("E18").Value = "  +3.73%  "
Set-TextValue $ws.Range("D19") "7.36"
$ws.Range("E19").Value = "  -2.03%  "
$ws.Range("E20").Value = "  -0.27%  "
Set-TextValue $ws.Range("D21") "10.86"
$ws.Range("E21").Value = "  -2.03%  "
Set-TextValue $ws.Range("D22") "468.19"
$ws.Range("E22").Value = "  -5.31%  "
Set-TextValue $ws.Range("D23") "0.738"
$ws.Range("E23").Value = "  +1.00%  "
Set-TextValue $ws.Range("D24") "0.0000160"
$ws.Range("E24").Value = "  -4.57%  "
Set-TextValue $ws.Range("D25") "83.46"
$ws.Range("E25").Value = "  -1.68%  "
Set-TextValue $ws.Range("D26") "2.23"
$ws.Range("E26").Value = "  -0.51%  "
Set-TextValue $ws.Range("D27") "12.11"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E29").Value = "  -2.20%  "
Set-TextValue $ws.Range("D30") "2.96"
$ws.Range("E30").Value = "  +0.23%  "
Set-TextValue $ws.Range("D31") "4.039.56"
$ws.Range("E31").Value = "  -0.92%  "
Set-TextValue $ws.Range("D32") "7.75"
$ws.Range("E32").Value = "  -0.08%  "
Set-TextValue $ws.Range("D33") "2.31"
$ws.Range("E33").Value = "  -3.02%  "
Set-TextValue $ws.Range("D34") "31.29"
$ws.Range("E34").Value = "  -2.29%  "
Set-TextValue $ws.Range("D35") "9.40"
$ws.Range("E35").Value = "  -0.35%  "
Set-TextValue $ws.Range("D36") "3.859.51"
$ws.Range("E36").Value = "  -0.63%  "
Set-TextValue $ws.Range("D37") "0.104"
$ws.Range("E37").Value = "  -2.61%  "
$ws.Range("E38").Value = "  +12.97%  "
$ws.Range("E39").Value = "  -0.92%  "
Set-TextValue $ws.Range("D40") "0.140"
$ws.Range("E40").Value = "  -0.19%  "
Set-TextValue $ws.Range("D41") "5.91"
$ws.Range("E41").Value = "  -1.04%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws.Range("D43") "0.312"
$ws.Range("E43").Value = "  -1.93%  "
$ws.Range("B44").Value = "FLOKI"
$ws.Range("C44").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
Set-TextValue $ws.Range("D44") "0.000302"
$ws.Range("E44").Value = "  +11.43%  "
Set-TextValue $ws.Range("D45") "1.98"
$ws.Range("E45").Value = "  -0.57%  "
Set-TextValue $ws.Range("D46") "423.04"
$ws.Range("E46").Value = "  -2.71%  "
$ws.Range("E47").Value = "  -0.01%  "
Set-TextValue $ws.Range("D48") "8.61"
$ws.Range("E48").Value = "  +0.07%  "
Set-TextValue $ws.Range("D49") "47.21"
$ws.Range("E49").Value = "  -1.84%  "
Set-TextValue $ws.Range("D50") "27.83"
$ws.Range("E50").Value = "  +6.88%  "
Set-TextValue $ws.Range("D51") "143.54"
$ws.Range("E51").Value = "  +0.08%  "

Write-Host "Applied crypto price updates"